$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Regime Atual"
$ws.Range("B2").Value = 0.6195539201671987
$ws.Range("C2").Value = 0.1396814823355647
$ws.Range("D2").Value = 0.5504218036392498
$ws.Range("E2").Value = 0.2474934668802131
$ws.Range("F2").Value = 314.466778929297

# Row 3 - "Nova Proposta"
$ws.Range("B3").Value = 0.6173387331698661
$ws.Range("C3").Value = 0.140375533103393
$ws.Range("D3").Value = 0.5476545770454402
$ws.Range("E3").Value = 0.2419947385555356
$ws.Range("F3").Value = 346.9600114092114
$ws.Range("G3").Value = 32.49323247991435

# Row 4 - "Nova c/ Aliq. Máxima"
$ws.Range("B4").Value = 0.616592356188896
$ws.Range("C4").Value = 0.1406526291330148
$ws.Range("D4").Value = 0.5467616641713284
$ws.Range("E4").Value = 0.2403697917822596
$ws.Range("F4").Value = 360.024117147753
$ws.Range("G4").Value = 45.55733821845598
